# Weekly data refresh: insert a new price record as row 88 (date 2023-01-31 /
# serial 44957) on "Vega Monumental Concepción" / Arándano (blue), shifting
# every existing row from 88 downward by one position (old row 88 becomes
# row 89, ..., old row 135 becomes row 136).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 88..135 down one slot, opening up a blank row 88.
$ws.Rows.Item(88).Insert()

# Populate the newly opened row 88 with the new weekly record.
$ws.Cells.Item(88, 1).Value  = 11
$ws.Cells.Item(88, 2).Value  = "Vega Monumental Concepción"
$ws.Cells.Item(88, 3).Value  = "Bíobío"
$ws.Cells.Item(88, 4).Value  = 44957
$ws.Cells.Item(88, 5).Value  = 8
$ws.Cells.Item(88, 6).Value  = "Fruta"
$ws.Cells.Item(88, 7).Value  = 100101
$ws.Cells.Item(88, 8).Value  = "Berries"
$ws.Cells.Item(88, 9).Value  = 100101001
$ws.Cells.Item(88, 10).Value = "Arándano (blue)"
$ws.Cells.Item(88, 11).Value = "Sin especificar"
$ws.Cells.Item(88, 12).Value = "Primera"
$ws.Cells.Item(88, 13).Value = 200
$ws.Cells.Item(88, 14).Value = 3000
$ws.Cells.Item(88, 15).Value = 3500
$ws.Cells.Item(88, 16).Value = 3250
$ws.Cells.Item(88, 17).Value = "$/bandeja 2 kilos"
$ws.Cells.Item(88, 18).Value = "Región de Ñuble"
$ws.Cells.Item(88, 19).Value = 1625
$ws.Cells.Item(88, 20).Value = 2
